$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes ---
# "Ultra-son arrière : 0x005" gains a " ??" suffix
$ws.Range("B18").Value = "Ultra-son arrière : 0x005 ??"

# The BUS CAN summary sentence is reworded
$ws.Range("B4").Value = "BUS CAN : Traitement différent en fonction des ids associés à l'envoie + récupération DATA"

# --- Remove the stray numeric 0 that used to sit next to "GPS : 0x004" ---
$ws.Range("C17").ClearContents()

# --- "ids :" label becomes bold ---
$ws.Range("B13").Font.Bold = $true

# --- Make the "Second temps :" lead-in of B32 bold, leaving the rest normal ---
$cell = $ws.Range("B32")
$fullLen = $cell.Value2.Length
$cell.Characters(1, 14).Font.Bold = $true
$cell.Characters(15, $fullLen - 14).Font.Bold = $false

# --- View changes: zoom level and active selection ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("E11").Select() | Out-Null
